$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header column C from "image" to "imageFile"
$ws.Range("C1").Value = "imageFile"

# Strip the "stim/" folder prefix from the image file names in column C (rows 2-9)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("stim/")) {
        $cell.Value = $val.ToString().Substring(5)
    }
}

# Move the active selection to C1
$ws.Range("C1").Select()

$wb.Save()
